# "Results from R script" — append two new observation rows (111, 112)
# to the bottom of the single data sheet, matching the existing
# date/volume/high/low/open/close/ticker layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date/time number format from the last existing data row
# (A110) onto the two new date cells first, so the new A111/A112 cells
# reuse the workbook's existing "yyyy-mm-dd hh:mm:ss" style record
# instead of minting a new (duplicate) one.
$ws.Range("A110").Copy()
$ws.Range("A111:A112").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
    @{ Row = 111; Date = 45474.2916666667; Volume = 0;   Price = 1.91999995708466; Ticker = "KK.MI" },
    @{ Row = 112; Date = 45475.3314236111; Volume = 900; Price = 1.91999995708466; Ticker = "KK.MI" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Date      # date
    $ws.Cells.Item($row, 2).Value = $r.Volume    # volume
    $ws.Cells.Item($row, 3).Value = $r.Price     # high
    $ws.Cells.Item($row, 4).Value = $r.Price     # low
    $ws.Cells.Item($row, 5).Value = $r.Price     # open
    $ws.Cells.Item($row, 6).Value = $r.Price     # close

    # adj_close: the value is numeric-looking but must be stored as
    # text (shared string), matching the rest of column G. Force text
    # via a "@" number format, then strip the format back off so the
    # cell keeps the plain default style (no explicit s="...") while
    # remaining a text cell.
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 7).Value = [string]$r.Price
    $ws.Cells.Item($row, 7).ClearFormats()

    $ws.Cells.Item($row, 8).Value = $r.Ticker    # ticker
}
